$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8395
$ws1.Range("F3").Value = 36597
$ws1.Range("G3").Value = 0
$ws1.Range("F12").Value = 662
$ws1.Range("F13").Value = 492
$ws1.Range("F17").Value = 449
$ws1.Range("F18").Value = 435
$ws1.Range("F22").Value = 2439
$ws1.Range("F30").Value = 1123

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 8395
$ws4.Range("F5").Value = 36597
$ws4.Range("G5").Value = 0
$ws4.Range("F18").Value = 662
$ws4.Range("F19").Value = 492
$ws4.Range("F28").Value = 449
$ws4.Range("F29").Value = 435
$ws4.Range("F33").Value = 2439
$ws4.Range("F42").Value = 1123
